$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Tên khách hàng" in H1, matching style of existing headers
$ws.Range("H1").Value = "Tên khách hàng"
$ws.Range("H1").HorizontalAlignment = $ws.Range("G1").HorizontalAlignment

# Set width for the new column H
$ws.Columns.Item(8).ColumnWidth = 15.3

# Update the active selection as in the saved view
$ws.Range("I6").Select()
